$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fruta / Hortaliza price data: each data row (2-24) is
# reassigned to a different day's record (dates reshuffled); row 23 unaffected.

# Row 2
$ws.Range("D2").Value = 44412
$ws.Range("J2").Value = 240
$ws.Range("N2").Value = '$/caja 40 unidades'
$ws.Range("P2").Value = 538
$ws.Range("Q2").Value = 40
# Row 3
$ws.Range("D3").Value = 44419
$ws.Range("H3").Value = 'Symphony'
$ws.Range("K3").Value = 21000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 21500
$ws.Range("N3").Value = '$/caja 50 unidades'
$ws.Range("P3").Value = 430
$ws.Range("Q3").Value = 50
# Row 4
$ws.Range("D4").Value = 44468
$ws.Range("H4").Value = 'Argentina(o)'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("N4").Value = '$/caja 50 unidades'
$ws.Range("P4").Value = 350
$ws.Range("Q4").Value = 50
# Row 5
$ws.Range("D5").Value = 44391
$ws.Range("J5").Value = 140
# Row 6
$ws.Range("D6").Value = 44405
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 21000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21500
$ws.Range("N6").Value = '$/caja 40 unidades'
$ws.Range("P6").Value = 538
$ws.Range("Q6").Value = 40
# Row 7
$ws.Range("D7").Value = 44377
$ws.Range("H7").Value = 'Madrigal'
$ws.Range("J7").Value = 150
$ws.Range("M7").Value = 20333
$ws.Range("P7").Value = 508
# Row 8
$ws.Range("D8").Value = 44377
$ws.Range("H8").Value = 'Symphony'
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21500
$ws.Range("P8").Value = 538
# Row 9
$ws.Range("D9").Value = 44384
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 21000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 21500
$ws.Range("P9").Value = 538
# Row 10
$ws.Range("D10").Value = 44384
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19333
$ws.Range("N10").Value = '$/caja 50 unidades'
$ws.Range("P10").Value = 387
$ws.Range("Q10").Value = 50
# Row 11
$ws.Range("D11").Value = 44384
$ws.Range("H11").Value = 'Symphony'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 21000
$ws.Range("M11").Value = 20400
$ws.Range("P11").Value = 510
# Row 12
$ws.Range("D12").Value = 44435
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 19000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19500
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("P12").Value = 488
$ws.Range("Q12").Value = 40
# Row 13
$ws.Range("D13").Value = 44398
$ws.Range("J13").Value = 170
$ws.Range("K13").Value = 21000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 21500
$ws.Range("P13").Value = 538
# Row 14
$ws.Range("D14").Value = 44433
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19500
$ws.Range("P14").Value = 488
# Row 15
$ws.Range("D15").Value = 44363
# Row 16
$ws.Range("D16").Value = 44356
$ws.Range("H16").Value = 'Argentina(o)'
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 19000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19500
$ws.Range("N16").Value = '$/caja 50 unidades'
$ws.Range("P16").Value = 390
$ws.Range("Q16").Value = 50
# Row 17
$ws.Range("D17").Value = 44482
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("P17").Value = 362
# Row 18
$ws.Range("D18").Value = 44426
$ws.Range("H18").Value = 'Madrigal'
$ws.Range("J18").Value = 150
$ws.Range("N18").Value = '$/caja 40 unidades'
$ws.Range("P18").Value = 488
$ws.Range("Q18").Value = 40
# Row 19
$ws.Range("D19").Value = 44483
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("P19").Value = 362
# Row 20
$ws.Range("D20").Value = 44489
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("P20").Value = 338
# Row 21
$ws.Range("D21").Value = 44370
$ws.Range("H21").Value = 'Argentina(o)'
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 21000
$ws.Range("M21").Value = 20429
$ws.Range("N21").Value = '$/caja 50 unidades'
$ws.Range("P21").Value = 409
$ws.Range("Q21").Value = 50
# Row 22
$ws.Range("D22").Value = 44370
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 22000
$ws.Range("L22").Value = 23000
$ws.Range("M22").Value = 22500
$ws.Range("P22").Value = 562
# Row 24
$ws.Range("D24").Value = 44160
$ws.Range("H24").Value = 'Madrigal'
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("N24").Value = '$/caja 40 unidades'
$ws.Range("P24").Value = 362
$ws.Range("Q24").Value = 40
